# Add a new "04-ago" column (AP) to the right of the existing "03-ago"
# column (AO), shifting the previous AO figures into AP and writing the
# freshly-reported figures into AO.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column AP (42nd column), mirroring the existing "03-ago"
# header style/format already applied to row 1.
$ws.Cells.Item(1, 42).Value = "04-ago"

# Updated figures for "03-ago" (column AO) and the new figures that land in
# the newly appended "04-ago" column (column AP), keyed by data row.
$aoValues = @{
    2  = 15
    3  = 16
    4  = 12
    5  = 15
    6  = 9
    7  = 18
    8  = 15
    9  = 17
    10 = 14
    11 = 13
}
$apValues = @{
    2  = 16
    3  = 14
    4  = 11
    5  = 17
    6  = 8
    7  = 18
    8  = 16
    9  = 16
    10 = 13
    11 = 12
}

# Match column AP's number format / alignment to column AO (numeric,
# centered) before writing data, same as every other data column.
$ws.Range("AP2:AP11").HorizontalAlignment = $ws.Range("AO2").HorizontalAlignment()
$ws.Range("AP2:AP11").NumberFormat = $ws.Range("AO2").NumberFormat()

foreach ($r in 2..11) {
    $ws.Cells.Item($r, 41).Value = $aoValues[$r]
    $ws.Cells.Item($r, 42).Value = $apValues[$r]
}

# Match the saved selection from the edited workbook.
$ws.Range("AS9").Select()
